# Planeacion cronograma proyecto - status + comments update
# (mirrors the author's "Add files via upload" commit: several WBS tasks
#  moved to "Complete" and received closing comments in the Comentarios
#  column.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example - Project Plan Template")

# --- 1. "Estado" (column D) -> "Complete" for the tasks that wrapped up ---
$completedRows = @(31, 36, 40, 42, 43, 44, 46, 47, 48, 50, 51, 52, 54, 55, 56)
foreach ($r in $completedRows) {
    $ws.Range("D$r").Value = "Complete"
}

# --- 2. New closing comments in "Comentarios" (column I) ---
$ws.Range("I46").Value = "Se identificaron varios errores menores que afectaban el rendimiento; fueron corregidos"
$ws.Range("I47").Value = "El proceso tomó más tiempo de lo previsto debido a incompatibilidades entre módulos y se lograron mejoras en la precisión del modelo tras ajustes en los parámetros."
$ws.Range("I48").Value = "Se redujo significativamente el tiempo de carga en dispositivos de gama media."
$ws.Range("I50").Value = "Se requirió actualizar librerías del servidor para compatibilidad con la última versión del backend."
$ws.Range("I51").Value = "Algunos usuarios reportaron dificultades menores en el proceso de instalación"
$ws.Range("I52").Value = "La mayoría de usuarios expresó satisfacción general, especialmente en la velocidad y usabilidad."
$ws.Range("I54").Value = "Se agregaron diagramas actualizados según los últimos cambios de arquitectura."
$ws.Range("I55").Value = "El informe integra de manera clara los avances, riesgos y decisiones tomadas."
$ws.Range("I56").Value = "Se ajustó el contenido para hacerlo más comprensible para stakeholders no técnicos."

# --- 3. Leave the view focused on the last-edited row, like the saved file ---
[void]$ws.Range("I52").Select()
